$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: the VIN in D5 was re-entered as a plain number (Excel stores it as a
# double, which only keeps ~15 significant digits of precision).
$ws.Range("D5").Value = 12345678901234570

# Row 6: new booking for the same customer as row 4 ("Ikki maru"), VIN entered
# as a plain number again (same underlying double as row 2's VIN).
$ws.Range("A6").Value = 631886740
$ws.Range("B6").Value = "Ikki maru"
$ws.Range("C6").Value = 992907510905
$ws.Range("D6").Value = 12345678912345680
$ws.Range("E6").Value = "Toyota Camry MILLION"
$ws.Range("F6").Value = "Регулярное обслуживание"
$ws.Range("G6").Value = "31/08/2025"
$ws.Range("H6").Value = "2025-08-22 14:44:41"
$ws.Range("I6").Value = "15:00"
$ws.Range("J6").Value = "-"

# Row 7: another new booking for the same customer, this time the VIN is kept
# as text. Build it as a formula returning the digit string, then paste back
# as a value-only copy so the cell ends up holding a plain text literal
# (no quote-prefix style attached), same as the rest of the sheet.
$ws.Range("A7").Value = 631886740
$ws.Range("B7").Value = "Ikki maru"
$ws.Range("C7").Value = 992907510905
$ws.Range("D7").Formula = '="12345678912345678"'
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "Toyota Land Cruiser"
$ws.Range("F7").Value = "Регулярное обслуживание"
$ws.Range("G7").Value = "23/08/2025"
$ws.Range("H7").Value = "2025-08-22 16:50:58"
$ws.Range("I7").Value = "15:00"
$ws.Range("J7").Value = "-"
